$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1053
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 1053
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 1053
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -1279
$ws.Range("H12").Value = 200
$ws.Range("I12").Value = 200
$ws.Range("K12").Value = 200
$ws.Range("M12").Value = -30
$ws.Range("H62").Value = 3897
$ws.Range("I62").Value = 3515.625
$ws.Range("K62").Value = 3515.625
$ws.Range("M62").Value = -2891.625
$ws.Range("H65").Value = 3897
$ws.Range("I65").Value = 3515.625
$ws.Range("K65").Value = 17578.125
$ws.Range("M65").Value = -14458.125
$ws.Range("H92").Value = 126011
$ws.Range("I92").Value = 126011
$ws.Range("K92").Value = 126011
$ws.Range("M92").Value = -124763
$ws.Range("H98").Value = 1601.8334
$ws.Range("I98").Value = 1601.8334
$ws.Range("K98").Value = 1601.8334
$ws.Range("M98").Value = -103.8334
$ws.Range("H101").Value = 763.1667
$ws.Range("I101").Value = 763.1667
$ws.Range("K101").Value = 2289.5001
$ws.Range("M101").Value = -667.5001000000002
$ws.Range("H103").Value = 206.77777
$ws.Range("I103").Value = 233.28572
$ws.Range("J103").Value = 114
$ws.Range("K103").Value = 699.85716
$ws.Range("L103").Value = 342
$ws.Range("M103").Value = -113.85716
$ws.Range("N103").Value = -1514
$ws.Range("H122").Value = 1601.8334
$ws.Range("I122").Value = 1601.8334
$ws.Range("K122").Value = 4805.5002
$ws.Range("M122").Value = -2355.5002
$ws.Range("H131").Value = 1896.3334
$ws.Range("I131").Value = 1528.75
$ws.Range("K131").Value = 4586.25
$ws.Range("M131").Value = 453.75
$ws.Range("H137").Value = 2025.5
$ws.Range("I137").Value = 2168
$ws.Range("J137").Value = 1940
$ws.Range("K137").Value = 6504
$ws.Range("L137").Value = 5820
$ws.Range("M137").Value = -3954
$ws.Range("N137").Value = -10920
$ws.Range("H138").Value = 4796.459
$ws.Range("I138").Value = 3247.303
$ws.Range("J138").Value = 6622.25
$ws.Range("K138").Value = 9741.909
$ws.Range("L138").Value = 19866.75
$ws.Range("M138").Value = -4601.909
$ws.Range("N138").Value = -30146.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 31720.715
$ws.Range("I55").Value = 25016
$ws.Range("J55").Value = 36749.25
$ws.Range("K55").Value = 25016
$ws.Range("L55").Value = 36749.25
$ws.Range("M55").Value = -24701
$ws.Range("N55").Value = -37379.25
$ws.Range("H132").Value = 1405.3334
$ws.Range("I132").Value = 1163.68
$ws.Range("J132").Value = 1954.5454
$ws.Range("K132").Value = 3491.04
$ws.Range("L132").Value = 5863.6362
$ws.Range("M132").Value = -961.04
$ws.Range("N132").Value = -10923.6362

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3496.2856
$ws.Range("I86").Value = 2949.818
$ws.Range("J86").Value = 5500
$ws.Range("K86").Value = 2949.818
$ws.Range("L86").Value = 5500
$ws.Range("M86").Value = -1826.818
$ws.Range("N86").Value = -7746
$ws.Range("H89").Value = 3496.2856
$ws.Range("I89").Value = 2949.818
$ws.Range("J89").Value = 5500
$ws.Range("K89").Value = 14749.09
$ws.Range("L89").Value = 27500
$ws.Range("M89").Value = -9133.09
$ws.Range("N89").Value = -38732
$ws.Range("H94").Value = 1539.6
$ws.Range("I94").Value = 1562
$ws.Range("K94").Value = 1562
$ws.Range("M94").Value = -1111
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H131").Value = 80000
$ws.Range("J131").Value = 80000
$ws.Range("L131").Value = 80000
$ws.Range("N131").Value = -90080

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H31").Value = 4047.1785
$ws.Range("I31").Value = 3526.2173
$ws.Range("J31").Value = 6443.6
$ws.Range("K31").Value = 3526.2173
$ws.Range("L31").Value = 6443.6
$ws.Range("M31").Value = -3231.2173
$ws.Range("N31").Value = -7033.6
$ws.Range("H34").Value = 4047.1785
$ws.Range("I34").Value = 3526.2173
$ws.Range("J34").Value = 6443.6
$ws.Range("K34").Value = 3526.2173
$ws.Range("L34").Value = 6443.6
$ws.Range("M34").Value = -3324.2173
$ws.Range("N34").Value = -6847.6
$ws.Range("H58").Value = 2428.6924
$ws.Range("I58").Value = 1228.8889
$ws.Range("K58").Value = 1228.8889
$ws.Range("M58").Value = -1025.8889
$ws.Range("H86").Value = 8499.333000000001
$ws.Range("J86").Value = 10249
$ws.Range("L86").Value = 10249
$ws.Range("N86").Value = -12495
$ws.Range("H89").Value = 8499.333000000001
$ws.Range("J89").Value = 10249
$ws.Range("L89").Value = 51245
$ws.Range("N89").Value = -62477
$ws.Range("H107").Value = 573.9
$ws.Range("I107").Value = 573.9
$ws.Range("K107").Value = 573.9
$ws.Range("M107").Value = 1346.1
$ws.Range("H136").Value = 2428.6924
$ws.Range("I136").Value = 1228.8889
$ws.Range("K136").Value = 3686.6667
$ws.Range("M136").Value = -1136.6667

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1585.875
$ws.Range("J5").Value = 1829.4
$ws.Range("L5").Value = 5488.200000000001
$ws.Range("N5").Value = -5712.200000000001
$ws.Range("H7").Value = 16666780
$ws.Range("I7").Value = 20000116
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 60000348
$ws.Range("L7").Value = 300
$ws.Range("M7").Value = -60000236
$ws.Range("N7").Value = -524
$ws.Range("H10").Value = 125.25
$ws.Range("I10").Value = 50
$ws.Range("J10").Value = 150.33333
$ws.Range("K10").Value = 150
$ws.Range("L10").Value = 450.99999
$ws.Range("M10").Value = -11
$ws.Range("N10").Value = -728.99999
$ws.Range("H26").Value = 272.66666
$ws.Range("I26").Value = 210.57143
$ws.Range("J26").Value = 490
$ws.Range("K26").Value = 631.71429
$ws.Range("L26").Value = 1470
$ws.Range("M26").Value = -343.71429
$ws.Range("N26").Value = -2046
$ws.Range("H34").Value = 1462.5385
$ws.Range("J34").Value = 2719
$ws.Range("L34").Value = 8157
$ws.Range("N34").Value = -8325
$ws.Range("H36").Value = 1487.4
$ws.Range("I36").Value = 1486
$ws.Range("J36").Value = 1500
$ws.Range("K36").Value = 4458
$ws.Range("L36").Value = 4500
$ws.Range("M36").Value = -4289
$ws.Range("N36").Value = -4838
$ws.Range("H39").Value = 9500
$ws.Range("I39").Value = 7666.6665
$ws.Range("K39").Value = 22999.9995
$ws.Range("M39").Value = -22705.9995
$ws.Range("H55").Value = 169523.33
$ws.Range("I55").Value = 252410
$ws.Range("J55").Value = 3750
$ws.Range("K55").Value = 757230
$ws.Range("L55").Value = 11250
$ws.Range("M55").Value = -757053
$ws.Range("N55").Value = -11604
$ws.Range("H107").Value = 1800.3334
$ws.Range("I107").Value = 1500
$ws.Range("J107").Value = 1950.5
$ws.Range("K107").Value = 4500
$ws.Range("L107").Value = 5851.5
$ws.Range("M107").Value = -2580
$ws.Range("N107").Value = -9691.5
$ws.Range("H122").Value = 308.375
$ws.Range("I122").Value = 442.33334
$ws.Range("K122").Value = 3981.00006
$ws.Range("M122").Value = -1531.00006
$ws.Range("H131").Value = 1250.6842
$ws.Range("I131").Value = 676
$ws.Range("K131").Value = 2028
$ws.Range("M131").Value = 3012
$ws.Range("H135").Value = 1585.875
$ws.Range("J135").Value = 1829.4
$ws.Range("L135").Value = 16464.6
$ws.Range("N135").Value = -21534.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 20000
$ws.Range("J49").Value = 20000
$ws.Range("L49").Value = 20000
$ws.Range("N49").Value = -20368
$ws.Range("H113").Value = 1766.909
$ws.Range("I113").Value = 929.625
$ws.Range("J113").Value = 3999.6667
$ws.Range("K113").Value = 929.625
$ws.Range("L113").Value = 3999.6667
$ws.Range("M113").Value = 1240.375
$ws.Range("N113").Value = -8339.6667
$ws.Range("H132").Value = 3349
$ws.Range("I132").Value = 2987.889
$ws.Range("K132").Value = 8963.667000000001
$ws.Range("M132").Value = -6433.667000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 6672666.5
$ws.Range("I2").Value = 20000000
$ws.Range("J2").Value = 9000
$ws.Range("K2").Value = 20000000
$ws.Range("L2").Value = 9000
$ws.Range("M2").Value = -19999888
$ws.Range("N2").Value = -9224
$ws.Range("H7").Value = 4346.25
$ws.Range("I7").Value = 4131.6665
$ws.Range("J7").Value = 4990
$ws.Range("K7").Value = 4131.6665
$ws.Range("L7").Value = 4990
$ws.Range("M7").Value = -4019.6665
$ws.Range("N7").Value = -5214
$ws.Range("H36").Value = 79888
$ws.Range("J36").Value = 79888
$ws.Range("L36").Value = 79888
$ws.Range("N36").Value = -81012
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H126").Value = 4346.25
$ws.Range("I126").Value = 4131.6665
$ws.Range("J126").Value = 4990
$ws.Range("K126").Value = 12394.9995
$ws.Range("L126").Value = 14970
$ws.Range("M126").Value = -9924.999500000002
$ws.Range("N126").Value = -19910

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 47000
$ws.Range("J123").Value = 47000
$ws.Range("L123").Value = 47000
$ws.Range("N123").Value = -56800
$ws.Range("H132").Value = 27919.344
$ws.Range("I132").Value = 41399.434
$ws.Range("J132").Value = 2082.5
$ws.Range("K132").Value = 124198.302
$ws.Range("L132").Value = 6247.5
$ws.Range("M132").Value = -121668.302
$ws.Range("N132").Value = -11307.5
